$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O4").Value = "Internal Assignment"
$ws.Range("O4").Font.Bold = $true
$ws.Range("O4").Font.Size = 12
$ws.Range("O4").Font.Name = "Calibri"
$ws.Range("O4").Font.Color = 0
Write-Host "done"
